# EPBDS-10604 Support Java 15
#
# On the "Rules" sheet, the throwNPE() example is expanded from a single
# trigger line into two lines describing the new (Java 9+) idiom:
#   before:  Integer.decode(null);
#   after :  Object arg = null;
#            java.util.Objects.requireNonNull(arg);
#
# This pushes every row below it down by one, so a blank row is inserted
# first (which naturally shifts the rest of the sheet and the sheet
# dimension down by one row), then the new/edited cell values are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Rules" sheet

# Insert a blank row at 23: old rows 23/24/25 become 24/25/26.
$ws.Rows.Item(23).Insert()

# The new blank row (B23) should look like the "code line" row above it
# (B22, still holding the old "Integer.decode(null);" text/style at this
# point), so copy that formatting down first.
$ws.Cells.Item(22, 2).Copy()
$ws.Cells.Item(23, 2).PasteSpecial(-4122)   # xlPasteFormats

# Fill in the text. B23 ("java.util.Objects.requireNonNull(arg);") is set
# before B22 ("Object arg = null;") so that the new shared strings are
# created in the same order as the target workbook.
$ws.Cells.Item(23, 2).Value2 = "java.util.Objects.requireNonNull(arg);"
$ws.Cells.Item(22, 2).Value2 = "Object arg = null;"
